# Generate Report for Handback
# Applies handback status / datetime / target+handback file links to the
# zh-cn and de-de localization-status sheets.

$wb = $excel.ActiveWorkbook

$mdUrlBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6c771129a269ce61dc8192a22e9b0b2dfe18f277/e2e/"
$file1 = "1675235d-1b88-4182-9bf0-a59470c29b50"
$file2 = "1f31a2b9-3103-4681-bab6-8d5dd01bb270"

$statusText = "Handed back: in sync with en-US"
$hyperlinkColor = 15570276   # BGR packed value of RGB(0x64,0x95,0xED) == style "HyperLink" font color

function Set-HandbackSheet($sheetName, $langSuffix, $hashFile1, $hashFile2, $handbackDateTime) {
    $ws = $wb.Worksheets.Item($sheetName)

    # --- Row 2 (file 1) ---
    $ws.Range("C2").Value = $statusText
    $ws.Range("I2").Value = "$file1.md"
    $ws.Range("I2").Font.Color = $hyperlinkColor
    $ws.Range("I2").Font.Underline = 2
    $ws.Range("J2").Value = "$file1.$hashFile1.$langSuffix.xlf"
    $ws.Range("K2").Value = $handbackDateTime

    # --- Row 3 (file 2) ---
    $ws.Range("C3").Value = $statusText
    $ws.Range("I3").Value = "$file2.md"
    $ws.Range("I3").Font.Color = $hyperlinkColor
    $ws.Range("I3").Font.Underline = 2
    $ws.Range("J3").Value = "$file2.$hashFile2.$langSuffix.xlf"
    $ws.Range("K3").Value = $handbackDateTime

    # --- Hyperlinks: rebuild so relationship ids line up the way Excel emits
    #     them when (re)writing hyperlinks for A2, I2, A3, I3 in that order.
    $ws.Hyperlinks.Delete()
    $ws.Hyperlinks.Add($ws.Range("A2"), "$mdUrlBase$file1.md", "", "", "$file1.md")
    $ws.Hyperlinks.Add($ws.Range("I2"), "$mdUrlBase$file1.md", "", "", "$file1.md")
    $ws.Hyperlinks.Add($ws.Range("A3"), "$mdUrlBase$file2.md", "", "", "$file2.md")
    $ws.Hyperlinks.Add($ws.Range("I3"), "$mdUrlBase$file2.md", "", "", "$file2.md")

    # --- Column widths (status / target-file / handback-file columns widened) ---
    $ws.Columns.Item(3).ColumnWidth = 29.166666666666668
    $ws.Columns.Item(9).ColumnWidth = 39.166666666666664
    $ws.Columns.Item(10).ColumnWidth = 39.166666666666664
}

Set-HandbackSheet "zh-cn" "zh-cn" "7bb48c34d267ca3973095f9ef595201c6600be55" "8de6312f368ab2c81725ca80d36b9d79db7471bd" "2016-11-09 06:33:38"
Set-HandbackSheet "de-de" "de-de" "7bb48c34d267ca3973095f9ef595201c6600be55" "8de6312f368ab2c81725ca80d36b9d79db7471bd" "2016-11-09 06:33:56"

# --- Overview sheet: the zh-cn / de-de columns are now wider to match ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668
